$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.251.32'
$ws.Range('E2').Value = '  -1.09%  '

$ws.Range('D3').Value = '3.316.64'
$ws.Range('E3').Value = '  +1.37%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').Value = '''186.84'
$ws.Range('E5').Value = '  +1.63%  '

$ws.Range('D6').Value = '''576.63'
$ws.Range('E6').Value = '  -0.76%  '

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').Value = '''0.603'
$ws.Range('E8').Value = '  -0.11%  '

$ws.Range('E9').Value = '  -0.50%  '

$ws.Range('D10').Value = '''6.64'
$ws.Range('E10').Value = '  +0.94%  '

$ws.Range('D11').Value = '''0.411'
$ws.Range('E11').Value = '  +0.49%  '

$ws.Range('D12').Value = '3.891.66'
$ws.Range('E12').Value = '  +1.47%  '

$ws.Range('E13').Value = '  -0.46%  '

$ws.Range('D14').Value = '''27.42'
$ws.Range('E14').Value = '  -0.07%  '

$ws.Range('D15').Value = '67.454.84'
$ws.Range('E15').Value = '  -0.83%  '

$ws.Range('E16').Value = '  -0.71%  '

$ws.Range('D17').Value = '3.325.41'
$ws.Range('E17').Value = '  +1.02%  '

$ws.Range('D18').Value = '''444.59'
$ws.Range('E18').Value = '  +8.93%  '

$ws.Range('E19').Value = '  +0.15%  '

$ws.Range('D20').Value = '''13.53'
$ws.Range('E20').Value = '  +0.67%  '

$ws.Range('D21').Value = '''7.78'
$ws.Range('E21').Value = '  +3.07%  '

$ws.Range('D22').Value = '''73.66'
$ws.Range('E22').Value = '  +3.57%  '

$ws.Range('D23').Value = '''1.00'
$ws.Range('E23').Value = '  +0.13%  '

$ws.Range('B24').Value = 'WrappedeETH'
$ws.Range('C24').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D24').Value = '3.449.84'
$ws.Range('E24').Value = '  +1.14%  '

$ws.Range('B25').Value = 'Polygon'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D25').Value = '''0.514'
$ws.Range('E25').Value = '  +0.99%  '

$ws.Range('E26').Value = '  +1.31%  '

$ws.Range('E27').Value = '  +0.27%  '

$ws.Range('D28').Value = '''9.07'
$ws.Range('E28').Value = '  -4.22%  '

$ws.Range('E29').Value = '  -0.11%  '

$ws.Range('E30').Value = '  +1.62%  '

$ws.Range('D31').Value = '''22.88'
$ws.Range('E31').Value = '  +0.73%  '

$ws.Range('E32').Value = '  -2.34%  '

$ws.Range('E33').Value = '  -0.06%  '

$ws.Range('E34').Value = '  -0.86%  '

$ws.Range('D35').Value = '''6.79'
$ws.Range('E35').Value = '  -1.73%  '

$ws.Range('D36').Value = '''1.52'
$ws.Range('E36').Value = '  +4.82%  '

$ws.Range('D37').Value = '''162.74'
$ws.Range('E37').Value = '  -0.93%  '

$ws.Range('D38').Value = '''27.44'
$ws.Range('E38').Value = '  +0.79%  '

$ws.Range('E39').Value = '  -2.37%  '

$ws.Range('E40').Value = '  -1.35%  '

$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '2.773.60'
$ws.Range('E41').Value = '  +3.73%  '

$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = '''4.47'
$ws.Range('E42').Value = '  -0.64%  '

$ws.Range('D43').Value = '''6.32'
$ws.Range('E43').Value = '  -0.67%  '

$ws.Range('D44').Value = '''24.93'
$ws.Range('E44').Value = '  +1.29%  '

$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').Value = '''40.15'
$ws.Range('E45').Value = '  -1.67%  '

$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').Value = '''0.0671'
$ws.Range('E46').Value = '  -1.07%  '

$ws.Range('E47').Value = '  -0.81%  '

$ws.Range('D48').Value = '''328.03'
$ws.Range('E48').Value = '  -1.94%  '

$ws.Range('D49').Value = '''0.0274'
$ws.Range('E49').Value = '  -0.19%  '

$ws.Range('D50').Value = '''0.991'
$ws.Range('E50').Value = '  +2.13%  '

$ws.Range('E51').Value = '  -0.76%  '
